# Apply updated cryptocurrency price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.508.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "'2.109.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'334.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").Value = "'0.5244"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.79%  "

$ws.Range("D8").Value = "'0.4525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.01%  "

$ws.Range("D9").Value = "'53.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.56%  "

$ws.Range("D10").Value = "'0.08999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "

$ws.Range("D11").Value = "'1.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "'24.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.01%  "

$ws.Range("D13").Value = "'2.101.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").Value = "'6.785"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "'7.820"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "'96.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").Value = "'0.06627"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").Value = "'19.31"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "'6.307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "'30.557.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("D24").Value = "'12.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("D25").Value = "'2.342"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("D26").Value = "'2.350.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").Value = "'22.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").Value = "'2.582"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").Value = "'163.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").Value = "'132.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").Value = "'1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("D32").Value = "'0.1075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("E33").Value = "  +6.36%  "

$ws.Range("D34").Value = "'6.175"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("D35").Value = "'3.942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").Value = "'10.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.84%  "

$ws.Range("D37").Value = "'0.02580"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("D38").Value = "'0.06835"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "'5.547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").Value = "'12.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").Value = "'0.2293"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").Value = "'0.6933"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").Value = "'1.246"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").Value = "'2.401"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.74%  "

$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("D46").Value = "'0.6424"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("D47").Value = "'14.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").Value = "'3.659"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("D50").Value = "'1.220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.31%  "

$ws.Range("D51").Value = "'83.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
